$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.213.55'
$ws.Range("E2").Value = '  +1.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.907.09'
$ws.Range("E3").Value = '  +2.03%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.44'
$ws.Range("E5").Value = '  +1.23%  '

$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5249'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3780'
$ws.Range("E8").Value = '  +3.19%  '

$ws.Range("E9").Value = '  +1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.25'
$ws.Range("E10").Value = '  +2.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8989'
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("E12").Value = '  +2.51%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.905.31'
$ws.Range("E13").Value = '  +1.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '95.21'
$ws.Range("E14").Value = '  +0.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.270'
$ws.Range("E15").Value = '  +0.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008655'
$ws.Range("E17").Value = '  +1.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.52'
$ws.Range("E18").Value = '  +2.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.280.51'
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.085'
$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.153.59'
$ws.Range("E22").Value = '  +1.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.64'
$ws.Range("E23").Value = '  +2.93%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.450'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.326'
$ws.Range("E25").Value = '  +11.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.69'
$ws.Range("E26").Value = '  -1.80%  '

$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.739'
$ws.Range("E27").Value = '  -1.95%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.14'
$ws.Range("E28").Value = '  +1.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.88'
$ws.Range("E29").Value = '  +1.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.973'
$ws.Range("E30").Value = '  +5.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.822'
$ws.Range("E31").Value = '  +2.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09238'
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8213'
$ws.Range("E33").Value = '  +9.95%  '

$ws.Range("E34").Value = '  +0.47%  '

$ws.Range("E35").Value = '  +7.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.986'
$ws.Range("E36").Value = '  +0.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.310'
$ws.Range("E37").Value = '  +2.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.605'
$ws.Range("E38").Value = '  +2.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5675'
$ws.Range("E39").Value = '  +1.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01991'
$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("E41").Value = '  +0.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.011'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.644'
$ws.Range("E43").Value = '  +0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.33'
$ws.Range("E44").Value = '  +3.06%  '

$ws.Range("E45").Value = '  +2.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4845'
$ws.Range("E46").Value = '  +1.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.25'
$ws.Range("E47").Value = '  +1.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.626'
$ws.Range("E49").Value = '  +4.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.61'
$ws.Range("E50").Value = '  +1.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.85'
$ws.Range("E51").Value = '  +1.34%  '
